$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats / xlPasteValues codes used below
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --- Row 6: ara / FR / 'ajnabiun / TRUE -------------------------------
$a6 = $ws.Cells.Item(6, 1)
$a6.Value2 = "ara"

$b6 = $ws.Cells.Item(6, 2)
$b6.Value2 = "FR"

# Leading apostrophe must be preserved literally (not treated as Excel's
# "quote prefix" marker), so build the value through a throwaway formula
# and then collapse it back down to a plain value.
$c6 = $ws.Cells.Item(6, 3)
$c6.Formula = '="''ajnabiun"'
$c6.Copy()
$c6.PasteSpecial($xlPasteValues)
$c6.WrapText = $true

$d6 = $ws.Cells.Item(6, 4)
$srcD = $ws.Cells.Item(2, 4)
$srcD.Copy()
$d6.PasteSpecial($xlPasteValues)

# --- Row 7: ara / NFR / ghayr 'ajnabiin / TRUE -------------------------
$a7 = $ws.Cells.Item(7, 1)
$a7.Value2 = "ara"

$b7 = $ws.Cells.Item(7, 2)
$b7.Value2 = "NFR"

$c7 = $ws.Cells.Item(7, 3)
$c7.Value2 = "ghayr 'ajnabiin"
$c6.Copy()
$c7.PasteSpecial($xlPasteFormats)

$d7 = $ws.Cells.Item(7, 4)
$srcD.Copy()
$d7.PasteSpecial($xlPasteValues)

# --- Column C width ------------------------------------------------------
# Target stored width is 18.85; the engine quantises ColumnWidth (character
# units) to steps of 1/6 on save, and 18 is the input that lands closest
# (stored width 18.8333..).
$ws.Columns.Item(3).ColumnWidth = 18

# --- Row heights (wrapped cells need the taller auto row height) -------
$ws.Rows.Item(6).RowHeight = 14.9
$ws.Rows.Item(7).RowHeight = 14.9

# --- Selection -------------------------------------------------------
$ws.Range("D8").Select()
